# Columns J:L ("b.se" / "ci" / "p.val") were laid out one position too far
# left (the p.val column's own cell was missing), so the b.se and ci
# values had slid into J/K and the p-value text ended up duplicated into
# L instead of having its own J cell - this affected the header row too.
#
# Fix: for every row (header + all data rows), rotate J,K,L one step to
# the right: new J = old L (the p-value label/text), new K = old J
# (b.se), new L = old K (ci). Column M is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 15; $r++) {
    $oldJ = $ws.Cells.Item($r, 10).Value2   # b.se column content
    $oldK = $ws.Cells.Item($r, 11).Value2   # ci column content
    $oldL = $ws.Cells.Item($r, 12).Value2   # p.val column content

    $ws.Cells.Item($r, 10).Value = $oldL   # J -> p.val
    $ws.Cells.Item($r, 11).Value = $oldJ   # K -> b.se (old J)
    $ws.Cells.Item($r, 12).Value = $oldK   # L -> ci (old K)
}
